$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.312.47"
$ws.Range("E2").Value = "'  +1.75%  "
$ws.Range("D3").Value = "'1.840.97"
$ws.Range("E3").Value = "'  +0.66%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D5").Value = "'243.05"
$ws.Range("E5").Value = "'  -0.67%  "
$ws.Range("D7").Value = "'0.9998"
$ws.Range("E7").Value = "'  +0.07%  "
$ws.Range("E8").Value = "'  -0.33%  "
$ws.Range("D9").Value = "'0.07526"
$ws.Range("E9").Value = "'  -1.60%  "
$ws.Range("D10").Value = "'23.28"
$ws.Range("E10").Value = "'  +0.41%  "
$ws.Range("D11").Value = "'0.07655"
$ws.Range("E11").Value = "'  -1.68%  "
$ws.Range("D12").Value = "'1.834.88"
$ws.Range("E12").Value = "'  +0.35%  "
$ws.Range("D13").Value = "'5.085"
$ws.Range("E13").Value = "'  +0.10%  "
$ws.Range("D14").Value = "'0.6869"
$ws.Range("E14").Value = "'  +1.28%  "
$ws.Range("D15").Value = "'89.29"
$ws.Range("E15").Value = "'  -3.98%  "
$ws.Range("D16").Value = "'6.306"
$ws.Range("E16").Value = "'  -1.94%  "
$ws.Range("D17").Value = "'29.315.80"
$ws.Range("E17").Value = "'  +1.78%  "
$ws.Range("D18").Value = "'0.000008248"
$ws.Range("E18").Value = "'  +0.42%  "
$ws.Range("D19").Value = "'2.088.92"
$ws.Range("E19").Value = "'  +0.80%  "
$ws.Range("D20").Value = "'234.06"
$ws.Range("E20").Value = "'  -2.95%  "
$ws.Range("E21").Value = "'  -0.29%  "
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("D23").Value = "'7.480"
$ws.Range("E23").Value = "'  +0.60%  "
$ws.Range("D24").Value = "'0.9995"
$ws.Range("E24").Value = "'  -0.05%  "
$ws.Range("D25").Value = "'0.1460"
$ws.Range("E25").Value = "'  -2.06%  "
$ws.Range("D26").Value = "'160.13"
$ws.Range("E26").Value = "'  -0.58%  "
$ws.Range("E27").Value = "'  +1.51%  "
$ws.Range("D28").Value = "'18.10"
$ws.Range("E28").Value = "'  -0.48%  "
$ws.Range("E29").Value = "'  -1.16%  "
$ws.Range("D30").Value = "'4.229"
$ws.Range("E30").Value = "'  +0.20%  "
$ws.Range("D31").Value = "'4.139"
$ws.Range("E31").Value = "'  -0.42%  "
$ws.Range("D32").Value = "'1.203"
$ws.Range("E32").Value = "'  +1.28%  "
$ws.Range("D33").Value = "'0.05138"
$ws.Range("E33").Value = "'  +0.56%  "
$ws.Range("D34").Value = "'0.7731"
$ws.Range("E34").Value = "'  +0.06%  "
$ws.Range("D35").Value = "'1.848"
$ws.Range("E35").Value = "'  -0.18%  "
$ws.Range("D36").Value = "'1.139"
$ws.Range("E36").Value = "'  +0.32%  "
$ws.Range("D37").Value = "'2.671"
$ws.Range("E37").Value = "'  -0.78%  "
$ws.Range("D38").Value = "'1.291.76"
$ws.Range("E38").Value = "'  +2.35%  "
$ws.Range("D39").Value = "'0.01845"
$ws.Range("E39").Value = "'  -0.38%  "
$ws.Range("E40").Value = "'  +0.13%  "
$ws.Range("D41").Value = "'0.9445"
$ws.Range("E41").Value = "'  -1.42%  "
$ws.Range("D42").Value = "'105.68"
$ws.Range("E42").Value = "'  -0.97%  "
$ws.Range("D43").Value = "'0.9993"
$ws.Range("E43").Value = "'  -0.01%  "
$ws.Range("D44").Value = "'5.651"
$ws.Range("E44").Value = "'  -6.25%  "
$ws.Range("D45").Value = "'9.687"
$ws.Range("E45").Value = "'  +0.36%  "
$ws.Range("D46").Value = "'1.989.23"
$ws.Range("E46").Value = "'  +0.86%  "
$ws.Range("D47").Value = "'0.5200"
$ws.Range("E47").Value = "'  +0.77%  "
$ws.Range("D48").Value = "'1.774"
$ws.Range("E48").Value = "'  +1.73%  "
$ws.Range("D49").Value = "'0.00000000121"
$ws.Range("E49").Value = "'  +1.97%  "
$ws.Range("D50").Value = "'63.32"
$ws.Range("E50").Value = "'  -0.95%  "
$ws.Range("D51").Value = "'0.05922"
$ws.Range("E51").Value = "'  +0.69%  "
